# Auto-generated script applying the meteocat data refresh diff
# Commit: "Update automatic: dades i banners [2026-02-10 18:50]"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes a plain-text value into a cell without letting Excel's
# smart-entry logic reinterpret pure "NN%" strings as a numeric percentage.
# It forces Text entry, then restores original formatting (number format +
# style) by pasting formats back from an adjacent cell that already carries
# the correct (unmodified) style, so the cell's style index is unaffected.
function Set-TextValue($range, [string]$text, $formatSource) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $formatSource.Copy()
    $range.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

$ws.Range("E2").Value = "2026-02-10 18:48:28"
$ws.Range("I2").Value = "36.7 mm"
$ws.Range("E3").Value = "2026-02-10 18:48:30"
$ws.Range("I3").Value = "23.2 mm"
$ws.Range("E4").Value = "2026-02-10 18:48:33"
$ws.Range("J4").Value = "1004.0 hPa"
$ws.Range("E5").Value = "2026-02-10 18:48:35"
$ws.Range("I5").Value = "30.2 mm"
$ws.Range("E6").Value = "2026-02-10 18:48:40"
$ws.Range("J6").Value = "1004.3 hPa"
$ws.Range("E7").Value = "2026-02-10 18:48:43"
Set-TextValue $ws.Range("H7") "69%" $ws.Range("I7")
$ws.Range("J7").Value = "1004.7 hPa"
$ws.Range("O7").Value = "15.1 °C"
$ws.Range("E8").Value = "2026-02-10 18:48:48"
$ws.Range("J8").Value = "1004.6 hPa"
$ws.Range("O8").Value = "11.8 °C"
$ws.Range("E9").Value = "2026-02-10 18:48:50"
Set-TextValue $ws.Range("H9") "95%" $ws.Range("I9")
$ws.Range("E10").Value = "2026-02-10 18:48:52"
$ws.Range("O10").Value = "10.1 °C"
$ws.Range("E11").Value = "2026-02-10 18:48:55"
$ws.Range("I11").Value = "0.3 mm"
$ws.Range("O11").Value = "7.6 °C"
$ws.Range("E12").Value = "2026-02-10 18:48:57"
$ws.Range("O12").Value = "9.0 °C"
$ws.Range("E13").Value = "2026-02-10 18:48:59"
$ws.Range("I13").Value = "4.7 mm"
$ws.Range("J13").Value = "1006.8 hPa"
$ws.Range("L13").Value = "18.7 km/h - 212º 18:11 TU"
$ws.Range("O13").Value = "5.1 °C"
$ws.Range("E14").Value = "2026-02-10 18:49:02"
$ws.Range("O14").Value = "13.3 °C"
$ws.Range("E15").Value = "2026-02-10 18:49:04"
$ws.Range("E16").Value = "2026-02-10 18:49:07"
Set-TextValue $ws.Range("H16") "86%" $ws.Range("I16")
$ws.Range("I16").Value = "23.5 mm"
$ws.Range("O16").Value = "1.0 °C"
$ws.Range("E17").Value = "2026-02-10 18:49:09"
Set-TextValue $ws.Range("H17") "93%" $ws.Range("I17")
$ws.Range("O17").Value = "4.6 °C"
$ws.Range("E18").Value = "2026-02-10 18:49:12"
$ws.Range("J18").Value = "1004.4 hPa"
$ws.Range("O18").Value = "10.2 °C"
$ws.Range("E19").Value = "2026-02-10 18:49:14"
$ws.Range("O19").Value = "6.4 °C"
$ws.Range("E20").Value = "2026-02-10 18:49:16"
$ws.Range("I20").Value = "6.3 mm"
$ws.Range("M20").Value = "2.8 °C 18:28 TU"
$ws.Range("O20").Value = "0.3 °C"
$ws.Range("E21").Value = "2026-02-10 18:49:19"
$ws.Range("I21").Value = "5.6 mm"
$ws.Range("E22").Value = "2026-02-10 18:49:21"
$ws.Range("I22").Value = "8.8 mm"
$ws.Range("E23").Value = "2026-02-10 18:49:24"
$ws.Range("I23").Value = "23.8 mm"
$ws.Range("E24").Value = "2026-02-10 18:49:26"
$ws.Range("E25").Value = "2026-02-10 18:49:28"
$ws.Range("I25").Value = "14.6 mm"
$ws.Range("E26").Value = "2026-02-10 18:49:31"
$ws.Range("E27").Value = "2026-02-10 18:49:33"
$ws.Range("I27").Value = "3.0 mm"
$ws.Range("O27").Value = "1.2 °C"
$ws.Range("E28").Value = "2026-02-10 18:49:36"
$ws.Range("J28").Value = "1004.3 hPa"
$ws.Range("O28").Value = "8.9 °C"
$ws.Range("E29").Value = "2026-02-10 18:49:38"
$ws.Range("E30").Value = "2026-02-10 18:49:41"
Set-TextValue $ws.Range("H30") "96%" $ws.Range("I30")
$ws.Range("J30").Value = "1004.4 hPa"
$ws.Range("O30").Value = "9.2 °C"
$ws.Range("E31").Value = "2026-02-10 18:49:43"
$ws.Range("J31").Value = "1003.6 hPa"
$ws.Range("E32").Value = "2026-02-10 18:49:46"
Set-TextValue $ws.Range("H32") "91%" $ws.Range("I32")
$ws.Range("E33").Value = "2026-02-10 18:49:48"
$ws.Range("I33").Value = "9.6 mm"
$ws.Range("J33").Value = "1006.6 hPa"
$ws.Range("O33").Value = "4.1 °C"
$ws.Range("E34").Value = "2026-02-10 18:49:51"
$ws.Range("I34").Value = "9.1 mm"
$ws.Range("O34").Value = "3.8 °C"
$ws.Range("E35").Value = "2026-02-10 18:49:53"
$ws.Range("E36").Value = "2026-02-10 18:49:56"
$ws.Range("J36").Value = "1004.5 hPa"
$ws.Range("E37").Value = "2026-02-10 18:49:58"
$ws.Range("E38").Value = "2026-02-10 18:50:00"
$ws.Range("O38").Value = "10.7 °C"
$ws.Range("E39").Value = "2026-02-10 18:50:03"
$ws.Range("I39").Value = "6.9 mm"
$ws.Range("O39").Value = "1.8 °C"
$ws.Range("E40").Value = "2026-02-10 18:50:05"
$ws.Range("I40").Value = "10.7 mm"
$ws.Range("J40").Value = "1007.1 hPa"
$ws.Range("E41").Value = "2026-02-10 18:50:08"
$ws.Range("J41").Value = "1004.8 hPa"
$ws.Range("O41").Value = "14.2 °C"
$ws.Range("E42").Value = "2026-02-10 18:50:10"
$ws.Range("O42").Value = "10.3 °C"
$ws.Range("E43").Value = "2026-02-10 18:50:12"
$ws.Range("O43").Value = "9.4 °C"
$ws.Range("E44").Value = "2026-02-10 18:50:15"
$ws.Range("I44").Value = "24.3 mm"
$ws.Range("E45").Value = "2026-02-10 18:50:17"
$ws.Range("I45").Value = "31.1 mm"
$ws.Range("E46").Value = "2026-02-10 18:50:20"
